$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 435, shifting existing rows 435:451 down to 436:452.
$ws.Rows.Item(435).Insert()

# Populate the newly inserted row 435 with the new data record.
$ws.Cells.Item(435, 1).Value = 10
$ws.Cells.Item(435, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(435, 3).Value = "La Araucanía"
$ws.Cells.Item(435, 4).Value = 45075
$ws.Cells.Item(435, 5).Value = 9
$ws.Cells.Item(435, 6).Value = 100112001
$ws.Cells.Item(435, 7).Value = "Berenjena"
$ws.Cells.Item(435, 8).Value = "Sin especificar"
$ws.Cells.Item(435, 9).Value = "Primera"
$ws.Cells.Item(435, 10).Value = 80
$ws.Cells.Item(435, 11).Value = 12000
$ws.Cells.Item(435, 12).Value = 12000
$ws.Cells.Item(435, 13).Value = 12000
$ws.Cells.Item(435, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(435, 15).Value = "Región del Maule"
$ws.Cells.Item(435, 16).Value = 300
$ws.Cells.Item(435, 17).Value = 40
$ws.Cells.Item(435, 18).Value = "Hortaliza"
